{"js": "// Fix the wrong due date: \"HW 4, Due March 4\" -> \"HW 4, Due March 3\".\n// The document's \"_GoBack\" bookmark also needs to move from the very\n// first paragraph (where it sat before the edit) to right after the\n// corrected \"HW 4, Due March 3\" run (where Word leaves it after the\n// author's last text edit).\n\nconst doc = context.document;\nconst body = doc.body;\n\n// 1) Find the due-date heading and fix the day.\nconst dateResults = body.search(\"HW 4, Due March 4\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  const dateRange = dateResults.items[0];\n  dateRange.insertText(\"HW 4, Due March 3\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Drop the stale \"_GoBack\" bookmark from the top of the document\n//    (harmless if it is not present).\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Re-locate the corrected text and drop a fresh \"_GoBack\" bookmark\n//    immediately after it, matching where Word leaves the mark after\n//    the last edit made to the document.\nconst fixedResults = body.search(\"HW 4, Due March 3\", { matchCase: true });\nfixedResults.load(\"items\");\nawait context.sync();\n\nif (fixedResults.items.length > 0) {\n  const fixedRange = fixedResults.items[0];\n  const endRange = fixedRange.getRange(Word.RangeLocation.end);\n  endRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Fix the wrong due date: \"HW 4, Due March 4\" -> \"HW 4, Due March 3\".\n# The document's \"_GoBack\" bookmark also needs to move from the very\n# first paragraph (where it sat before the edit) to right after the\n# corrected \"HW 4, Due March 3\" run (where Word leaves it after the\n# author's last text edit).\n\n$d = $word.ActiveDocument\n\n# 1) Locate the due-date heading run.\n$headingRange = $d.Content\n$find = $headingRange.Find\n$find.Text = \"HW 4, Due March 4\"\n$found = $find.Execute()\n\nif ($found) {\n    # 2) Re-point \"_GoBack\" to sit right after the heading text *before*\n    #    editing the text itself. Bookmarks.Add relocates the bookmark in\n    #    place if one with this name already exists (it does, at the top\n    #    of the document), so this both removes the old one and creates\n    #    the new one in a single call, and it keeps the run that follows\n    #    (the lone space run, then \"The manipulator \") from getting\n    #    collapsed together when the text edit below happens.\n    $bmPoint = $d.Range($headingRange.End, $headingRange.End)\n    $d.Bookmarks.Add(\"_GoBack\", $bmPoint)\n\n    # 3) Now correct the due date text in place.\n    $headingRange2 = $d.Content\n    $find2 = $headingRange2.Find\n    $find2.Text = \"HW 4, Due March 4\"\n    $found2 = $find2.Execute()\n    if ($found2) {\n        $headingRange2.Text = \"HW 4, Due March 3\"\n    }\n}\n"}
